$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Restyle a couple of existing rows (B156 and B163 pick up the "style 1" look) ---
$ws.Range("B157").Copy()
$ws.Range("B156").PasteSpecial(-4122)
$ws.Range("B156").PasteSpecial(-4122)
$ws.Range("B163").PasteSpecial(-4122)

# --- Append three new employee rows ---
$ws.Range("A164").Value = 5290
$ws.Range("C164").Value = "ef77c5f6-bf69-4c0a-973e-021a0f09c1a3"
$ws.Range("B164").Value = "PALOMA LUCIA DOS SANTOS"

$ws.Range("A165").Value = 5291
$ws.Range("C165").Value = "784cd1a3-c037-4151-94de-bba47591cc3d"
$ws.Range("B165").Value = "DAYANE ARAUJO JESUS"

$ws.Range("A166").Value = 5292
$ws.Range("B166").Value = "SILEDI MARIA ALVES DOS SANTOS"
$ws.Range("C166").Value = "4f4f7a56-ae5d-485e-a4f1-6d35905142f4"

$ws.Range("B164:B166").Select()
